$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the main title text.
$ws.Range("A2").Value = "Calendario reuniones de seguimiento"

# 2) The whole sheet shifts up by one row (there is a leading blank row
#    being removed at the very top of the workbook's used layout).
$ws.Rows.Item(1).Delete()

# --- From here on, row numbers refer to the POST-shift numbering
#     (i.e. what used to be row N is now row N-1). The second table
#     ("Fechas Importantes") used to occupy rows 14-19 with columns
#     A=Semana, B:C=Motivo (merged), D=Fecha; it is now at rows 13-18
#     with columns A=Motivo, B=Fecha, C=Semana (no merges) and column D
#     dropped entirely.

# 3) Shrink the section title merge from A:D to A:C.
$ws.Range("A13:D13").UnMerge()
$ws.Range("A13:C13").Merge()
$ws.Range("D13").Clear()

# 4) Rearrange the header row (14) and the four data rows (15-18):
#    new column A <- old column B (Motivo)
#    new column B <- old column D (Fecha)
#    new column C <- old column A (Semana)
#    old column D is dropped entirely (content + formatting).
For ($r = 14; $r -le 18; $r++) {
    $ws.Range("B$r`:C$r").UnMerge()
    $ws.Range("A$r").Copy($ws.Range("F$r"))
    $ws.Range("B$r").Copy($ws.Range("A$r"))
    $ws.Range("D$r").Copy($ws.Range("B$r"))
    $ws.Range("F$r").Copy($ws.Range("C$r"))
    $ws.Range("D$r").Clear()
    $ws.Range("F$r").Clear()
}

# 5) Re-apply the bold/bordered look to the new 3-column header row.
$ws.Range("A14:C14").Font.Bold = $true
$ws.Range("A14:C14").HorizontalAlignment = -4108
$ws.Range("A14:C14").Borders.Item(7).LineStyle = 1
$ws.Range("A14:C14").Borders.Item(8).LineStyle = 1
$ws.Range("A14:C14").Borders.Item(9).LineStyle = 1
$ws.Range("A14:C14").Borders.Item(10).LineStyle = 1

# 6) Give the merged section-title box (row 13) its outer border + bold text.
$ws.Range("A13:C13").Font.Bold = $true
$ws.Range("A13:C13").HorizontalAlignment = -4108
$ws.Range("A13:C13").Borders.Item(8).LineStyle = 1
$ws.Range("A13:C13").Borders.Item(9).LineStyle = 1
$ws.Range("A13").Borders.Item(7).LineStyle = 1
$ws.Range("C13").Borders.Item(10).LineStyle = 1

Write-Host "done"
